$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1582278481012658
$ws.Range("C2").Value = 0.6645569620253164
$ws.Range("P2").Value = 0.08860759493670886
$ws.Range("S2").Value = 0.08860759493670886
$ws.Range("C3").Value = 0.05454545454545454
$ws.Range("J3").Value = 0.01818181818181818
$ws.Range("P3").Value = 0.8
$ws.Range("S3").Value = 0.1272727272727273
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.03389830508474576
$ws.Range("D6").Value = 0.01129943502824859
$ws.Range("F6").Value = 0.04519774011299435
$ws.Range("J6").Value = 0.2203389830508475
$ws.Range("O6").Value = 0.03389830508474576
$ws.Range("Q6").Value = 0.2090395480225989
$ws.Range("R6").Value = 0.06779661016949153
$ws.Range("S6").Value = 0.3785310734463277
$ws.Range("B7").Value = 0.06779661016949153
$ws.Range("D7").Value = 0.01129943502824859
$ws.Range("F7").Value = 0.07344632768361582
$ws.Range("J7").Value = 0.1299435028248588
$ws.Range("O7").Value = 0.02824858757062147
$ws.Range("Q7").Value = 0.2598870056497175
$ws.Range("R7").Value = 0.04519774011299435
$ws.Range("S7").Value = 0.384180790960452
$ws.Range("B8").Value = 0.03696098562628337
$ws.Range("D8").Value = 0.01848049281314168
$ws.Range("F8").Value = 0.03080082135523614
$ws.Range("J8").Value = 0.1006160164271047
$ws.Range("O8").Value = 0.01026694045174538
$ws.Range("Q8").Value = 0.2443531827515401
$ws.Range("R8").Value = 0.07392197125256673
$ws.Range("S8").Value = 0.484599589322382
$ws.Range("B9").Value = 0.04810996563573883
$ws.Range("D9").Value = 0.003436426116838488
$ws.Range("E9").Value = 0.003436426116838488
$ws.Range("F9").Value = 0.04123711340206185
$ws.Range("J9").Value = 0.1099656357388316
$ws.Range("O9").Value = 0.01718213058419244
$ws.Range("Q9").Value = 0.2508591065292096
$ws.Range("R9").Value = 0.05498281786941581
$ws.Range("S9").Value = 0.4707903780068728
$ws.Range("B10").Value = 0.07064760302775441
$ws.Range("D10").Value = 0.01682085786375105
$ws.Range("E10").Value = 0.0008410428931875525
$ws.Range("F10").Value = 0.06980656013456686
$ws.Range("J10").Value = 0.1286795626576955
$ws.Range("O10").Value = 0.01177460050462574
$ws.Range("Q10").Value = 0.2514718250630782
$ws.Range("R10").Value = 0.07232968881412952
$ws.Range("S10").Value = 0.3776282590412111
$ws.Range("G11").Value = 0.1434977578475336
$ws.Range("J11").Value = 0.08968609865470852
$ws.Range("K11").Value = 0.1704035874439462
$ws.Range("L11").Value = 0.5874439461883408
$ws.Range("S11").Value = 0.008968609865470852
$ws.Range("G12").Value = 0.8161764705882353
$ws.Range("J12").Value = 0.1397058823529412
$ws.Range("L12").Value = 0.02941176470588235
$ws.Range("S12").Value = 0.01470588235294118
$ws.Range("G13").Value = 0.8222222222222222
$ws.Range("J13").Value = 0.1555555555555556
$ws.Range("S13").Value = 0.02222222222222222
$ws.Range("H15").Value = 0.1944444444444444
$ws.Range("I15").Value = 0.1157407407407407
$ws.Range("J15").Value = 0.3148148148148148
$ws.Range("K15").Value = 0.06944444444444445
$ws.Range("M15").Value = 0.03240740740740741
$ws.Range("O15").Value = 0.04166666666666666
$ws.Range("S15").Value = 0.2314814814814815
$ws.Range("F16").Value = 0.008333333333333333
$ws.Range("H16").Value = 0.1833333333333333
$ws.Range("I16").Value = 0.1
$ws.Range("J16").Value = 0.475
$ws.Range("M16").Value = 0.01666666666666667
$ws.Range("O16").Value = 0.05833333333333333
$ws.Range("S16").Value = 0.09166666666666666
$ws.Range("F17").Value = 0.008741258741258742
$ws.Range("H17").Value = 0.201048951048951
$ws.Range("I17").Value = 0.1328671328671329
$ws.Range("J17").Value = 0.4073426573426573
$ws.Range("K17").Value = 0.06643356643356643
$ws.Range("M17").Value = 0.01923076923076923
$ws.Range("O17").Value = 0.06993006993006994
$ws.Range("S17").Value = 0.0944055944055944
$ws.Range("H18").Value = 0.2531645569620253
$ws.Range("I18").Value = 0.1518987341772152
$ws.Range("J18").Value = 0.3354430379746836
$ws.Range("K18").Value = 0.0379746835443038
$ws.Range("M18").Value = 0.0189873417721519
$ws.Range("O18").Value = 0.08860759493670886
$ws.Range("S18").Value = 0.1139240506329114
$ws.Range("F19").Value = 0.01857835218093699
$ws.Range("H19").Value = 0.2189014539579968
$ws.Range("I19").Value = 0.1268174474959612
$ws.Range("J19").Value = 0.3610662358642973
$ws.Range("K19").Value = 0.09289176090468497
$ws.Range("M19").Value = 0.01777059773828756
$ws.Range("N19").Value = 0.003231017770597738
$ws.Range("O19").Value = 0.0630048465266559
$ws.Range("S19").Value = 0.09773828756058159
